$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D/E columns hold text like "26.892.04" / "  +1.84%  " that must stay text
# (NumberFormat "@" keeps Excel from silently recasting them as numbers).
# Apply once, up front, over the whole data range so every touched cell
# gets a single consistent text style.
$ws.Range("D2:E51").NumberFormat = "@"

# Updated price/volume data for rows 2-33 (values refreshed by the scraper).

$data = @(
    @{ Row = 2;  D = "26.892.04";    E = "  +1.84%  " },
    @{ Row = 3;  D = "1.727.02";     E = "  +0.31%  " },
    @{ Row = 4;  D = "0.9974";       E = "  -0.26%  " },
    @{ Row = 5;  D = "241.98";       E = "  -0.43%  " },
    @{ Row = 6;  D = "0.9981";       E = "  -0.22%  " },
    @{ Row = 7;  D = "0.4897";       E = "  -0.20%  " },
    @{ Row = 8;  D = "0.2590";       E = "  -0.75%  " },
    @{ Row = 9;  D = "0.06214";      E = "  +0.43%  " },
    @{ Row = 10; D = "1.729.90";     E = "  +0.47%  " },
    @{ Row = 11; D = "15.99";        E = "  +3.36%  " },
    @{ Row = 12; D = "0.06903";      E = "  -1.57%  " },
    @{ Row = 13; D = "0.6090";       E = "  +1.76%  " },
    @{ Row = 14; E = "  -1.61%  " },
    @{ Row = 15; D = "77.23";        E = "  +0.06%  " },
    @{ Row = 16; D = "0.9983";       E = "  -0.20%  " },
    @{ Row = 17; D = "26.873.32";    E = "  +1.77%  " },
    @{ Row = 18; D = "0.9974";       E = "  -0.28%  " },
    @{ Row = 19; D = "0.000007184";  E = "  +0.76%  " },
    @{ Row = 20; D = "11.45";        E = "  +0.99%  " },
    @{ Row = 21; D = "1.952.80";     E = "  +0.46%  " },
    @{ Row = 22; E = "  -1.01%  " },
    @{ Row = 23; D = "8.575";        E = "  +0.03%  " },
    @{ Row = 24; D = "5.101";        E = "  -0.97%  " },
    @{ Row = 25; D = "138.65";       E = "  +1.00%  " },
    @{ Row = 26; D = "15.31";        E = "  +0.70%  " },
    @{ Row = 27; D = "1.794";        E = "  +5.39%  " },
    @{ Row = 28; E = "  -1.19%  " },
    @{ Row = 29; D = "106.07";       E = "  -0.81%  " },
    @{ Row = 30; D = "3.947";        E = "  +0.19%  " },
    @{ Row = 31; D = "0.07995";      E = "  +0.61%  " },
    @{ Row = 32; D = "3.688";        E = "  +0.61%  " },
    @{ Row = 33; D = "0.04533" }
)

foreach ($item in $data) {
    $r = $item.Row
    if ($item.ContainsKey("D")) {
        $ws.Range("D$r").Value = $item.D
    }
    if ($item.ContainsKey("E")) {
        $ws.Range("E$r").Value = $item.E
    }
}

# Rows 34-51: the coin list shifted up by one row (the "Frax" row was
# removed), and a new coin (NEARProtocol) was appended at the bottom (row 51).
$coins = @(
    @{ Row = 34; B = "HuobiToken";        C = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht";        D = "2.597";   E = "  -0.27%  " },
    @{ Row = 35; B = "ARBITRUM";          C = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb";             D = "1.008";   E = "  +1.54%  " },
    @{ Row = 36; B = "ImmutableX";        C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx";           D = "0.6255";  E = "  +0.38%  " },
    @{ Row = 37; B = "TrustWalletToken";  C = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt";     D = "0.9380";  E = "  +1.24%  " },
    @{ Row = 38; B = "RenderToken";       C = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr";     D = "2.056";   E = "  +5.77%  " },
    @{ Row = 39; B = "MXToken";           C = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx";            D = "2.457";   E = "  +2.75%  " },
    @{ Row = 40; B = "PaxDollar";         C = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp";           D = "0.9980";  E = "  -0.19%  " },
    @{ Row = 41; B = "VeChain";           C = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet";          D = "0.01503"; E = "  +1.44%  " },
    @{ Row = 42; B = "FraxShare";         C = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs";            D = "5.653";   E = "  +5.95%  " },
    @{ Row = 43; B = "Quant";             C = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt";            D = "99.46";   E = "  -0.47%  " },
    @{ Row = 44; B = "TheSandbox";        C = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand";          D = "0.3853";  E = "  +0.50%  " },
    @{ Row = 45; B = "Aptos";             C = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt";                D = "6.886";   E = "  +2.72%  " },
    @{ Row = 46; B = "Algorand";          C = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo";        D = "0.1163";  E = "  +0.16%  " },
    @{ Row = 47; B = "Cronos";            C = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro";            D = "0.05395"; E = "  +0.58%  " },
    @{ Row = 48; B = "EnergySwap";        C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens";           D = "7.907";   E = "  +2.96%  " },
    @{ Row = 49; B = "Elrond";            C = "https://coinranking.com/coin/omwkOTglq+elrond-egld";              D = "30.18";   E = "  +0.38%  " },
    @{ Row = 50; B = "Aave";              C = "https://coinranking.com/coin/ixgUfzmLR+aave-aave";                D = "51.66";   E = "  +1.74%  " },
    @{ Row = 51; B = "NEARProtocol";      C = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near";        D = "1.236";   E = "  +0.14%  " }
)

foreach ($item in $coins) {
    $r = $item.Row
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
}
